# EnterpriseJavaTimeLog.xlsx edit script
# Implements: Hibernate Search prep rows, NPE troubleshooting/fix rows, Friday note,
# plus compaction of the trailing "Issues/Loose Ends" notes block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) First clear every row from 45 downward so nothing old lingers, then we
#    rebuild the tail of the sheet from scratch in its new, shorter layout.
# ---------------------------------------------------------------------------
$ws.Range("A45:D70").Clear()

# ---------------------------------------------------------------------------
# 1) Row 42/43: shared-string text unchanged, only the underlying index shifted
#    (an earlier string got removed) -- re-assert the same text so nothing is lost.
# ---------------------------------------------------------------------------
$ws.Range("D42").Value = "Indie Project: Brought Survey Dao Tests to 5/5 passing; Noted programming decisions that are needed for searching surveys.  Revised tables (and other files as needed) to structure the roles table as neede for authentication.  Created data for the application database.`nWeek 7: Started videos, following along in project."
$ws.Range("D43").Value = "Week 7/Project: worked on web.xml changes and package changes for authentication."

# ---------------------------------------------------------------------------
# 2) Row 44: hours revised upward, note text expanded, row grows taller.
# ---------------------------------------------------------------------------
$ws.Range("B44").Value = 5.5
$ws.Range("D44").Value = "Week 7: Completed readings (light treatment) and videos; added hibernate search to pom, made cfg changes, added annotations`nIndie Project: added separate jsp's for 403 errors and other errors`nProfessional Development: Researching Hibernate Search; started a ppt"
$ws.Rows.Item(44).RowHeight = 60

# ---------------------------------------------------------------------------
# 3) New rows 45 & 46 - two new time-log entries.
# ---------------------------------------------------------------------------
$ws.Range("A45").Value = 43541
$ws.Range("A45").NumberFormat = "d-mmm"
$ws.Range("B45").Value = 3
$ws.Range("D45").Value = "Indie Project: trying to figure out how to create index from existing database for Hibernate Search.  Troubleshooting NPE that I swear wasn't there before."
$ws.Range("D45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 30

$ws.Range("A46").Value = 43546
$ws.Range("A46").NumberFormat = "d-mmm"
$ws.Range("B46").Value = 4
$ws.Range("D46").Value = "Eliminated the NPE!!!"
$ws.Range("D46").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Row 47 stays a blank date-formatted row.
# ---------------------------------------------------------------------------
$ws.Range("A47").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------------
# 5) Row 48: blank date cell retained, height reverts to default, note text replaced.
# ---------------------------------------------------------------------------
$ws.Range("A48").NumberFormat = "d-mmm"
$ws.Range("D48").Value = "Friday 9pm - x…"
$ws.Range("D48").WrapText = $true

# ---------------------------------------------------------------------------
# 6) Rows 49/50 stay untouched blank date rows.
# ---------------------------------------------------------------------------
$ws.Range("A49").NumberFormat = "d-mmm"
$ws.Range("A50").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------------
# 7) Row 51: becomes a standalone note row (old "A" date cell removed),
#    carries the "TODO (nice to do)" note and a custom row height.
# ---------------------------------------------------------------------------
$ws.Range("D51").Value = "TODO (nice to do)  -  refactor tests to use .equals in UserDaoTest and StoryDaoTest"
$ws.Range("D51").WrapText = $true
$ws.Rows.Item(51).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 8) Rows 52/54 used to be blank date rows; in the new layout they don't
#    exist at all (already cleared above, nothing more to do).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 9) Remaining "Issues / Loose Ends" notes get compacted upward by a few rows
#    each, keeping their own text/format but moving to new row numbers.
# ---------------------------------------------------------------------------
$ws.Range("D53").Value = "Issues/Loose Ends:"
$ws.Range("D53").WrapText = $true

$ws.Range("D55").Value = "2019-03-01 06:11:44,732 Log4j2-TF-1-RollingFileManager-1 ERROR Error in post-rollover Delete when visiting C:\logs\income_experiences.log java.nio.file.FileSystemException: C:\logs\income_experiences.log: The process cannot access the file because it is being used by another process."
$ws.Range("D55").WrapText = $true
$ws.Rows.Item(55).RowHeight = 45

$ws.Range("D56").Value = "^^^ Recently unable to re-create"
$ws.Range("D56").WrapText = $true

$ws.Range("D58").Value = "NB: new password so files/setup need to change in prior repos"
$ws.Range("D58").WrapText = $true

$ws.Range("D60").Value = "Removed from UserDao comments -- save for a later time, differetn place"
$ws.Range("D60").WrapText = $true

$ws.Range("D61").Value = ""

$ws.Range("D63").Value = "Useful mysql tutorial: used for changing foreign key constraint:"
$ws.Range("D63").WrapText = $true

$ws.Range("D64").Value = "http://www.mysqltutorial.org/mysql-foreign-key/"
$ws.Range("D64").WrapText = $true

# ---------------------------------------------------------------------------
# 10) View state: scroll/selection roughly matches the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("D49").Select()
